$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A30").Value = 46000
$ws.Range("B30").Value = 62

$ws.Range("A30:B30").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
